# Apply the NoteTypeVS ValueSet update:
#  - bump Version to 0.1.15-beta
#  - bump Date to 2023-06-07T11:47:17-05:00
#  - remove the "83320-2 / Allergy and Immunology Adverse event note" row
#    from the "Include from LOINC" table (allergy notes are excluded by
#    business rule, not translated)

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B3").Value = "0.1.15-beta"
$wsMeta.Range("B8").Value = "2023-06-07T11:47:17-05:00"

$wsLoinc = $wb.Worksheets.Item("Include from LOINC")
$wsLoinc.Rows("2:2").Delete()
